# Add a "Shipment Type" column into the invoice export template.
# The new column is inserted before the existing "VoucherID" column (old
# column E), shifting all the following header/data/placeholder columns
# one slot to the right (old E..Y -> new F..Z).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; Excel copies the formatting (incl. width) of
# the column immediately to its left (D), same as a manual "Insert Column".
$ws.Columns("E:E").Insert()

# Match column D's width (the inserted column otherwise gets the default
# width since no explicit width carries over from the insert).
$ws.Columns("E:E").ColumnWidth = 29.15

# New header cell (row 1) and new placeholder cell (row 2).
$ws.Range("E1").Value = "Shipment Type"
$ws.Range("E2").Value = "{ShipmentType}"

# Give the two new cells their own distinct formatting (bold + centered +
# bordered for the header, bordered for the data row) so they don't keep
# sharing the style used by the rest of row 1 / row 2.
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").Font.Size = 12
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").Borders.LineStyle = 1

$ws.Range("E2").Font.Bold = $false
$ws.Range("E2").Font.Size = 12
$ws.Range("E2").Borders.LineStyle = 1

# Selection moves back to A2 (top-left of the frozen data pane).
$ws.Range("A2").Select()

# Page orientation set to portrait.
$ws.PageSetup.Orientation = 1
